$d = $word.ActiveDocument

# 1) Swap the placeholder ID text and drop the trailing space run by
#    replacing the whole "id + trailing space" span in one shot - Word
#    collapses the match back down to a single run.
[void]$d.Content.Find.Execute("**ID__AFFARS_5301_topic_29__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5301_90__ID**", 2)

# 2) First paragraph: bump the left indent from 120 twips (6pt) to 225
#    twips (11.25pt) and give it a paragraph border with 5-twip
#    (0.25pt) spacing on all four sides.
$p = $d.Paragraphs.Item(1)
$p.LeftIndent = 11.25

$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
